# Applies crypto price/volume update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value as literal text, even when it looks numeric
# (e.g. "1.210", "0.00001200"), without perturbing the cell style.
function Set-TextValue($range, [string]$value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "23.426.02"
$ws.Range("E2").Value = "  -1.50%  "

Set-TextValue $ws.Range("D3") "1.643.45"
$ws.Range("E3").Value = "  -0.77%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("E5").Value = "  +0.06%  "

Set-TextValue $ws.Range("D6") "298.65"
$ws.Range("E6").Value = "  -1.82%  "

Set-TextValue $ws.Range("D7") "0.3776"
$ws.Range("E7").Value = "  -1.34%  "

Set-TextValue $ws.Range("D8") "50.07"
$ws.Range("E8").Value = "  -1.98%  "

Set-TextValue $ws.Range("D9") "0.3520"
$ws.Range("E9").Value = "  -2.61%  "

Set-TextValue $ws.Range("D10") "0.08064"
$ws.Range("E10").Value = "  -1.96%  "

Set-TextValue $ws.Range("D11") "1.210"
$ws.Range("E11").Value = "  -3.87%  "

$ws.Range("E12").Value = "  +0.13%  "

Set-TextValue $ws.Range("D13") "22.03"
$ws.Range("E13").Value = "  -2.98%  "

Set-TextValue $ws.Range("D14") "6.370"
$ws.Range("E14").Value = "  -2.75%  "

Set-TextValue $ws.Range("D15") "7.300"
$ws.Range("E15").Value = "  -2.35%  "

Set-TextValue $ws.Range("D16") "0.00001200"
$ws.Range("E16").Value = "  -3.31%  "

Set-TextValue $ws.Range("D17") "1.638.22"
$ws.Range("E17").Value = "  +0.19%  "

Set-TextValue $ws.Range("D18") "96.69"
$ws.Range("E18").Value = "  -1.17%  "

Set-TextValue $ws.Range("D19") "0.06975"
$ws.Range("E19").Value = "  -0.06%  "

Set-TextValue $ws.Range("D20") "6.712"
$ws.Range("E20").Value = "  -1.05%  "

Set-TextValue $ws.Range("D21") "17.34"
$ws.Range("E21").Value = "  -2.46%  "

$ws.Range("E22").Value = "  -0.01%  "

Set-TextValue $ws.Range("D23") "12.36"
$ws.Range("E23").Value = "  -3.20%  "

Set-TextValue $ws.Range("D24") "23.438.93"
$ws.Range("E24").Value = "  -1.45%  "

Set-TextValue $ws.Range("D25") "2.480"
$ws.Range("E25").Value = "  -3.69%  "

Set-TextValue $ws.Range("D26") "2.889"
$ws.Range("E26").Value = "  -6.48%  "

Set-TextValue $ws.Range("D27") "20.82"
$ws.Range("E27").Value = "  -2.40%  "

Set-TextValue $ws.Range("D28") "153.06"
$ws.Range("E28").Value = "  +1.37%  "

Set-TextValue $ws.Range("D29") "5.205"
$ws.Range("E29").Value = "  -0.47%  "

Set-TextValue $ws.Range("D30") "132.15"
$ws.Range("E30").Value = "  -1.81%  "

Set-TextValue $ws.Range("D31") "1.818.62"
$ws.Range("E31").Value = "  -0.09%  "

Set-TextValue $ws.Range("D32") "6.880"
$ws.Range("E32").Value = "  -1.16%  "

Set-TextValue $ws.Range("D33") "2.137"
$ws.Range("E33").Value = "  -1.75%  "

Set-TextValue $ws.Range("D34") "11.53"
$ws.Range("E34").Value = "  -2.86%  "

Set-TextValue $ws.Range("D35") "0.9832"
$ws.Range("E35").Value = "  -9.19%  "

Set-TextValue $ws.Range("D36") "0.02706"
$ws.Range("E36").Value = "  -4.52%  "

$ws.Range("E37").Value = "  -1.22%  "

Set-TextValue $ws.Range("D38") "0.2435"
$ws.Range("E38").Value = "  -3.44%  "

Set-TextValue $ws.Range("D39") "5.906"
$ws.Range("E39").Value = "  -4.23%  "

Set-TextValue $ws.Range("D40") "0.06785"
$ws.Range("E40").Value = "  -5.39%  "

Set-TextValue $ws.Range("D41") "12.82"
$ws.Range("E41").Value = "  -2.57%  "

Set-TextValue $ws.Range("D42") "0.6846"
$ws.Range("E42").Value = "  -3.24%  "

$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D43") "1.290"
$ws.Range("E43").Value = "  -3.98%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D44") "15.59"
$ws.Range("E44").Value = "  -2.87%  "

$ws.Range("E45").Value = "  +0.15%  "

Set-TextValue $ws.Range("D46") "0.6327"
$ws.Range("E46").Value = "  -3.40%  "

Set-TextValue $ws.Range("D47") "2.247"
$ws.Range("E47").Value = "  -3.69%  "

Set-TextValue $ws.Range("D48") "3.898"
$ws.Range("E48").Value = "  -1.63%  "

Set-TextValue $ws.Range("D49") "0.07716"
$ws.Range("E49").Value = "  -3.24%  "

Set-TextValue $ws.Range("D50") "126.95"
$ws.Range("E50").Value = "  -1.31%  "

Set-TextValue $ws.Range("D51") "1.141"
$ws.Range("E51").Value = "  -4.39%  "
